# Prepare the Excel controller: add a new "06020" worksheet right after the
# existing "6020" sheet and seed it with the standard header row used for the
# line-editing template.

$wb = $excel.ActiveWorkbook

# Locate the existing "6020" sheet so the new sheet can be inserted after it.
$existingSheet = $wb.Worksheets.Item("6020")

# Add a new worksheet positioned immediately after "6020".
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $existingSheet)
$ws.Name = "06020"

# Seed the header row.
$ws.Range("A1").Value = "Omschrijving"
$ws.Range("B1").Value = "Inhoud"
$ws.Range("C1").Value = "Weergave"
$ws.Range("D1").Value = "Uitlijnen"
$ws.Range("E1").Value = "Regel verwijderen"
$ws.Range("F1").Value = "Regel template"
